# Update "想去人数" (want-to-go count) values on the 展览 and 全部类型 sheets
$wb = $excel.ActiveWorkbook

$sheetUpdates = @{
    "展览"    = @{ "F2" = 8606; "F6" = 1378; "F10" = 9378; "F12" = 99; "F15" = 357; "F16" = 6357; "F18" = 87; "F20" = 134 }
    "全部类型" = @{ "F2" = 8606; "F6" = 1378; "F12" = 9378; "F14" = 99; "F17" = 357; "F18" = 6357; "F20" = 87; "F22" = 134 }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $sheetUpdates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
